# This script finishes the "grouping" logic for the water-source-availability
# breakdown table. Row 4 (Rural_Al_Jiblah) and Row 5 (Urban_Al_Habelien) are
# re-computed against the new "C2_water_source_availability" grouping, and the
# previous (ungrouped) totals that used to live in row 4/5 are relocated into
# two brand-new rows (6 and 7) that are labelled with the
# "C2_water_source_availability" = no / yes categories. The old placeholder
# rows 15/16 (which only held those two labels) are removed since the data
# now lives in rows 6/7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 (Rural_Al_Jiblah) : updated values ----
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 867
$ws.Cells.Item(4, 10).Value = 0.10610696365194
$ws.Cells.Item(4, 11).Value = 6679
$ws.Cells.Item(4, 12).Value = 0.817403010647412
$ws.Cells.Item(4, 13).Value = 472
$ws.Cells.Item(4, 14).Value = 0.0577652674091299
$ws.Cells.Item(4, 17).Value = 153
$ws.Cells.Item(4, 18).Value = 0.0187247582915188
$ws.Cells.Item(4, 19).Value = 8171

# ---- Row 5 (Urban_Al_Habelien) : updated values ----
$ws.Cells.Item(5, 3).Value = 151
$ws.Cells.Item(5, 4).Value = 0.0145332050048123
$ws.Cells.Item(5, 5).Value = 545
$ws.Cells.Item(5, 6).Value = 0.0524542829643888
$ws.Cells.Item(5, 7).Value = 66
$ws.Cells.Item(5, 8).Value = 0.00635226179018287
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 49
$ws.Cells.Item(5, 12).Value = 0.00471607314725698
$ws.Cells.Item(5, 13).Value = 145
$ws.Cells.Item(5, 14).Value = 0.0139557266602502
$ws.Cells.Item(5, 16).Value = 0.026948989412897
$ws.Cells.Item(5, 17).Value = 9154
$ws.Cells.Item(5, 18).Value = 0.881039461020212
$ws.Cells.Item(5, 19).Value = 10390

# ---- Row 6 : new data row for C2_water_source_availability = "no" ----
$ws.Cells.Item(6, 1).Value = "C2_water_source_availability"
$ws.Cells.Item(6, 2).Value = "no"
$ws.Cells.Item(6, 3).Value = 66
$ws.Cells.Item(6, 4).Value = 0.00702351814408854
$ws.Cells.Item(6, 5).Value = 55
$ws.Cells.Item(6, 6).Value = 0.00585293178674045
$ws.Cells.Item(6, 7).Value = 66
$ws.Cells.Item(6, 8).Value = 0.00702351814408854
$ws.Cells.Item(6, 9).Value = 421
$ws.Cells.Item(6, 10).Value = 0.0448015324039587
$ws.Cells.Item(6, 11).Value = 6084
$ws.Cells.Item(6, 12).Value = 0.647440672555071
$ws.Cells.Item(6, 13).Value = 606
$ws.Cells.Item(6, 14).Value = 0.064488666595722
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 2099
$ws.Cells.Item(6, 18).Value = 0.223369160370331
$ws.Cells.Item(6, 19).Value = 9397
$ws.Cells.Item(6, 20).Value = 1

# ---- Row 7 : new data row for C2_water_source_availability = "yes" ----
$ws.Cells.Item(7, 2).Value = "yes"
$ws.Cells.Item(7, 3).Value = 85
$ws.Cells.Item(7, 4).Value = 0.00927542557835007
$ws.Cells.Item(7, 5).Value = 490
$ws.Cells.Item(7, 6).Value = 0.0534701003928416
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 446
$ws.Cells.Item(7, 10).Value = 0.0486687036228721
$ws.Cells.Item(7, 11).Value = 644
$ws.Cells.Item(7, 12).Value = 0.0702749890877346
$ws.Cells.Item(7, 13).Value = 11
$ws.Cells.Item(7, 14).Value = 0.00120034919249236
$ws.Cells.Item(7, 15).Value = 280
$ws.Cells.Item(7, 16).Value = 0.0305543430816237
$ws.Cells.Item(7, 17).Value = 7208
$ws.Cells.Item(7, 18).Value = 0.786556089044086
$ws.Cells.Item(7, 19).Value = 9164
$ws.Cells.Item(7, 20).Value = 1

# ---- Remove the old placeholder rows that only held the "no"/"yes" labels ----
# (row 16 first would also work since Delete() shifts rows up, but deleting
#  row 15 twice collapses both rows 15 and 16 cleanly)
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()

Write-Host "Applied grouping changes to rows 4-7 and removed old label rows 15-16"
